$d = $word.ActiveDocument

# Locate the "Introduction and Motivations" heading paragraph and bump its
# font size to 14pt (w:sz 28 half-points), applied to both the run and the
# paragraph mark (pPr/rPr), matching what Word does when the whole
# paragraph (incl. end-of-paragraph mark) is selected and resized.
$rng = $d.Content
$found = $rng.Find.Execute("Introduction and Motivations", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$para = $rng.Paragraphs(1)
$para.Range.Font.Size = 14
